# Applies corrected validation results to auditoria_validacion.xlsx
# Commit: "Se corrigen las validaciones"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: backtest_ranges
# Update pred_MAE / pred_RMSE / pred_MAPE and the matching
# pred_price_* columns (which mirror the same values) for rows 2-4.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("backtest_ranges")

$ws1.Range("F2").Value = 0.003735560000000006
$ws1.Range("G2").Value = 0.005003107714211245
$ws1.Range("H2").Value = 0.3388102804134973
$ws1.Range("L2").Value = 0.003735560000000006
$ws1.Range("M2").Value = 0.005003107714211245
$ws1.Range("N2").Value = 0.3388102804134973

$ws1.Range("F3").Value = 0.004470761978683473
$ws1.Range("G3").Value = 0.006090099392694067
$ws1.Range("H3").Value = 0.4055358688775248
$ws1.Range("L3").Value = 0.004470761978683473
$ws1.Range("M3").Value = 0.006090099392694067
$ws1.Range("N3").Value = 0.4055358688775248

$ws1.Range("F4").Value = 0.01580027341018982
$ws1.Range("G4").Value = 0.02099419234259793
$ws1.Range("H4").Value = 1.43649118887672
$ws1.Range("L4").Value = 0.01580027341018982
$ws1.Range("M4").Value = 0.02099419234259793
$ws1.Range("N4").Value = 1.43649118887672

# ---------------------------------------------------------------
# Sheet 2: signals_distribution
# Fill in the n / p_buy / p_hold / p_sell columns (C:F) for every
# data row - previously only the "n" column held a value for the
# buy/sell signal_col rows, and the hold rows were entirely empty.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("signals_distribution")

$ws2.Range("C2").Value = 500
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 100
$ws2.Range("F2").Value = 0

$ws2.Range("C3").Value = 500
$ws2.Range("D3").Value = 48.6
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 51.4

$ws2.Range("C4").Value = 500
$ws2.Range("D4").Value = 39.2
$ws2.Range("E4").Value = 22.4
$ws2.Range("F4").Value = 38.4

$ws2.Range("C5").Value = 500
$ws2.Range("D5").Value = 47.59999999999999
$ws2.Range("E5").Value = 0
$ws2.Range("F5").Value = 52.40000000000001

$ws2.Range("C6").Value = 500
$ws2.Range("D6").Value = 48.6
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = 51.4

$ws2.Range("C7").Value = 500
$ws2.Range("D7").Value = 33
$ws2.Range("E7").Value = 29.6
$ws2.Range("F7").Value = 37.4

$ws2.Range("C8").Value = 500
$ws2.Range("D8").Value = 37.6
$ws2.Range("E8").Value = 0
$ws2.Range("F8").Value = 62.4

$ws2.Range("C9").Value = 500
$ws2.Range("D9").Value = 48.6
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = 51.4

$ws2.Range("C10").Value = 500
$ws2.Range("D10").Value = 34.59999999999999
$ws2.Range("E10").Value = 5.2
$ws2.Range("F10").Value = 60.2

# ---------------------------------------------------------------
# Sheet 3: dm_tests
# Updated Diebold-Mariano statistics/p-values after the fix.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("dm_tests")

$ws3.Range("C2").Value = [double]"1.61394316183885E-16"

$ws3.Range("C4").Value = [double]"-1.670833306185265E-15"
$ws3.Range("D4").Value = 0.9999999999999987
